$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "cold space SOS – Scene01"

$ws.Range("F10").Select()
